$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A617").Value = 157
$ws.Range("B617").Value = "Output Properties"
$ws.Range("C617").Value = "Proprietà di uscita"

$ws.Range("B618").Value = "Audio"
$ws.Range("C618").Value = "Audio"

$ws.Range("B619").Value = "Video"
$ws.Range("C619").Value = "Video"

$ws.Range("B620").Value = "Disable frame rate limit"
$ws.Range("C620").Value = "Disattivare il limite della frequenza dei fotogrammi"

$ws.Range("B621").Value = "Disable sound output"
$ws.Range("C621").Value = "Disattivare l'uscita audio"

$ws.Range("B622").Value = "Capture before filtering"
$ws.Range("C622").Value = "Cattura prima del filtraggio"

$ws.Range("B623").Value = "Disable sound sync"
$ws.Range("C623").Value = "Disattivare la sincronizzazione del suono"

$ws.Range("B624").Value = "AVI output enabled"
$ws.Range("C624").Value = "Uscita AVI abilitata"

$ws.Range("B625").Value = "Ripper"
$ws.Range("C625").Value = "Squartatore"

$ws.Range("B626").Value = "Save screenshot"
$ws.Range("C626").Value = "Salva screenshot"

$ws.Range("B627").Value = "Pro Wizard 1.62"
$ws.Range("C627").Value = "Pro Wizard 1.62"

$ws.Range("B628").Value = "Sample ripper"
$ws.Range("C628").Value = "Ripper campione"

$ws.Range("B629").Value = "Take screenshot before filtering"
$ws.Range("C629").Value = "Scattare un'istantanea prima del filtraggio"

$ws.Range("B630").Value = "Autoclip screenshot"
$ws.Range("C630").Value = "Schermata Autoclip"

$ws.Range("B631").Value = "Re-recorder"
$ws.Range("C631").Value = "Ri-registratore"

$ws.Range("B632").Value = "Play recording"
$ws.Range("C632").Value = "Riproduzione della registrazione"

$ws.Range("B633").Value = "Automatic replay"
$ws.Range("C633").Value = "Riproduzione automatica"

$ws.Range("B634").Value = "Recording rate (seconds):"
$ws.Range("C634").Value = "Velocità di registrazione (secondi):"

$ws.Range("B635").Value = "Re-recording enabled"
$ws.Range("C635").Value = "Registrazione abilitata"

$ws.Range("B636").Value = "Save recording"
$ws.Range("C636").Value = "Salvare la registrazione"

$ws.Range("B637").Value = "Recording buffers:"
$ws.Range("C637").Value = "Buffer di registrazione:"

$ws.Range("A638").Value = 160
$ws.Range("B638").Value = "Filter Settings"
$ws.Range("C638").Value = "Impostazioni del filtro"

$ws.Range("B639").Value = "Reset to defaults"
$ws.Range("C639").Value = "Ripristino delle impostazioni predefinite"

$ws.Range("B640").Value = "Horiz. size:"
$ws.Range("C640").Value = "Dimensione orizzontale:"

$ws.Range("B641").Value = "Vert. size:"
$ws.Range("C641").Value = "Dimensione Vert:"

$ws.Range("B642").Value = "Horiz. position:"
$ws.Range("C642").Value = "Posizione orizzontale:"

$ws.Range("B643").Value = "Vert. position:"
$ws.Range("C643").Value = "Posizione verticale:"

$ws.Range("B644").Value = "Aspect Ratio Correction"
$ws.Range("C644").Value = "Correzione del rapporto d'aspetto"

$ws.Range("B645").Value = "Keep autoscale aspect"
$ws.Range("C645").Value = "Mantenere l'aspetto di autoscala"

$ws.Range("B646").Value = "Keep aspect ratio"
$ws.Range("C646").Value = "Mantenere il rapporto d'aspetto"

$ws.Range("B647").Value = "Extra Settings"
$ws.Range("C647").Value = "Impostazioni extra"

$ws.Range("B648").Value = "Presets"
$ws.Range("C648").Value = "Preimpostazioni"

$ws.Range("B649").Value = "Load"
$ws.Range("C649").Value = "Caricare"

$ws.Range("B650").Value = "Save"
$ws.Range("C650").Value = "Salvare"

$ws.Range("B651").Value = "Delete"
$ws.Range("C651").Value = "Cancellare"

$ws.Range("A652").Value = 163

# Update selection to match the post-edit cursor position
$ws.Range("A653").Select()